# Restore revision: change the "R30" rule's lower bound (cell C10 on the
# "Rules" sheet) from 18 to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C10").Value = 1
